# Edit: insert two new daily price records for "Pimiento" (Zafiro rojo / Zafiro verde)
# at Macroferia Regional de Talca, pushing the existing rows 254-322 down to 256-324.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 254 (old rows 254..322 become 256..324)
$ws.Rows("254:255").Insert()

# ---- New row 254: Zafiro rojo ----
$ws.Range("A254").Value = 5
$ws.Range("B254").Value = "Macroferia Regional de Talca"
$ws.Range("C254").Value = "Maule"
$ws.Range("D254").Value = 44463
$ws.Range("E254").Value = 7
$ws.Range("F254").Value = 100112002
$ws.Range("G254").Value = "Pimiento"
$ws.Range("H254").Value = "Zafiro rojo"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 200
$ws.Range("K254").Value = 40000
$ws.Range("L254").Value = 40000
$ws.Range("M254").Value = 40000
$ws.Range("N254").Value = "$/caja 15 kilos"
$ws.Range("O254").Value = "Región de Arica y Parinacota"
$ws.Range("P254").Value = 2667
$ws.Range("Q254").Value = 15
$ws.Range("R254").Value = "Hortaliza"

# ---- New row 255: Zafiro verde ----
$ws.Range("A255").Value = 5
$ws.Range("B255").Value = "Macroferia Regional de Talca"
$ws.Range("C255").Value = "Maule"
$ws.Range("D255").Value = 44463
$ws.Range("E255").Value = 7
$ws.Range("F255").Value = 100112002
$ws.Range("G255").Value = "Pimiento"
$ws.Range("H255").Value = "Zafiro verde"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 300
$ws.Range("K255").Value = 30000
$ws.Range("L255").Value = 30000
$ws.Range("M255").Value = 30000
$ws.Range("N255").Value = "$/caja 15 kilos"
$ws.Range("O255").Value = "Región de Arica y Parinacota"
$ws.Range("P255").Value = 2000
$ws.Range("Q255").Value = 15
$ws.Range("R255").Value = "Hortaliza"
